# Weekly update: a new price record for Camote (Vega Modelo de Temuco) is
# inserted at row 99, pushing the existing rows 99-116 down to rows 100-117.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 99 (shifts rows 99:116 down to 100:117).
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with this week's data.
$ws.Cells.Item(99, 1).Value = 10
$ws.Cells.Item(99, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(99, 3).Value = "La Araucanía"
$ws.Cells.Item(99, 4).Value = 44798
$ws.Cells.Item(99, 4).NumberFormat = $ws.Cells.Item(100, 4).NumberFormat
$ws.Cells.Item(99, 5).Value = 9
$ws.Cells.Item(99, 6).Value = 100114002
$ws.Cells.Item(99, 7).Value = "Camote"
$ws.Cells.Item(99, 8).Value = "Sin especificar"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 80
$ws.Cells.Item(99, 11).Value = 20000
$ws.Cells.Item(99, 12).Value = 20000
$ws.Cells.Item(99, 13).Value = 20000
$ws.Cells.Item(99, 14).Value = "`$/malla 20 kilos"
$ws.Cells.Item(99, 15).Value = "Perú"
$ws.Cells.Item(99, 16).Value = 1000
$ws.Cells.Item(99, 17).Value = 20
$ws.Cells.Item(99, 18).Value = "Hortaliza"
